$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# --- Add new rows demonstrating --skipby ---
# Row 11 (new shared string "--skipby(Key, 4)" introduced here first)
$ws.Range("A11").Value = "9"
$ws.Range("B11").Value = "css"
$ws.Range("C11").Value = "body > div.container-fluid > div > main > div:nth-child(36) > a.btn.btn-primary.btn-lg.active"
$ws.Range("D11").Value = "checkout"
$ws.Range("E11").Value = "--skipby(Key, 4)"
$ws.Range("F11").Value = "no_key"
$ws.Range("J11").Value = "fail"

# Row 12
$ws.Range("A12").Value = "10"
$ws.Range("B12").Value = "css"
$ws.Range("C12").Value = "none path, key is no, will jump"
$ws.Range("D12").Value = "checkout"
$ws.Range("E12").Value = "--skipby(Key, 4)"
$ws.Range("F12").Value = "no_key"
$ws.Range("J12").Value = "fail"

# Row 13 - index placeholder only
$ws.Range("A13").Value = "11"

# Row 14 - index placeholder only
$ws.Range("A14").Value = "12"

# Row 15 - index placeholder only (new shared string "13")
$ws.Range("A15").Value = "13"

# Row 16 - index placeholder only (new shared string "14")
$ws.Range("A16").Value = "14"

# Row 17 - final skip target row (new shared strings "15" and "--jumpto(Key, 0)")
$ws.Range("A17").Value = "15"
$ws.Range("B17").Value = "css"
$ws.Range("C17").Value = "none path, key is no, will jump"
$ws.Range("D17").Value = "checkout"
$ws.Range("E17").Value = "--jumpto(Key, 0)"
$ws.Range("F17").Value = "no_key"

# --- Update existing rows ---
# Row 6: change the "none path" example + its logic to the new skipby example
$ws.Range("E6").Value = "--skipby(No, 2)"
$ws.Range("C6").Value = "none path, will skip to seven"

# Row 10: jump target changed from 4 to 10
$ws.Range("E10").Value = "--jumpto(Key, 10)"

# --- Update selection to reflect the new active cell ---
$ws.Activate()
$ws.Range("E10").Select()
